# Apply the 3.3.1 RMI update to the Boolean Include Emissions from
# Imported Electricity (BIEfIE) control-settings workbook.

$wb = $excel.ActiveWorkbook

$aboutSheet = $wb.Worksheets.Item("About")
$biefieSheet = $wb.Worksheets.Item("BIEfIE")

# Stamp the "About" sheet with the last-edited date (2021-04-21) in C1,
# formatted as a date.
$aboutSheet.Range("C1").NumberFormat = "mm-dd-yy"
$aboutSheet.Range("C1").Value = "4/21/2021"

# Flip the BIEfIE boolean lever from 1 (include) to 0 (exclude).
$biefieSheet.Range("B2").Value = 0

# Reset BIEfIE's stale B3 selection back to the default top-left cell,
# then make "About" the selected / active sheet.
$biefieSheet.Activate()
$biefieSheet.Range("A1").Select()
$aboutSheet.Activate()
